# Refresh the cryptos price/volume table to match the latest scrape.
# Also reorders two coin pairs that swapped rank: row 35/36
# (WEMIXToken now above Hedera) and row 38/39 (Celestia now above Kaspa).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "39.237.40"
    "E2" = "  -1.68%  "
    "D3" = "2.187.98"
    "E3" = "  -6.65%  "
    "E4" = "  +0.09%  "
    "D5" = "295.95"
    "E5" = "  -3.98%  "
    "D6" = "81.44"
    "E6" = "  -3.07%  "
    "E7" = "  -4.01%  "
    "E8" = "  +0.06%  "
    "D9" = "0.465"
    "E9" = "  -3.80%  "
    "D10" = "0.0767"
    "E10" = "  -6.26%  "
    "D11" = "29.04"
    "E11" = "  -3.26%  "
    "D12" = "47.12"
    "E12" = "  -10.28%  "
    "E13" = "  -2.22%  "
    "E14" = "  -2.60%  "
    "D15" = "2.526.99"
    "E15" = "  -6.78%  "
    "D16" = "13.91"
    "E16" = "  -6.10%  "
    "D17" = "2.186.11"
    "E17" = "  -7.40%  "
    "D18" = "0.707"
    "D19" = "39.109.49"
    "E19" = "  -2.17%  "
    "D20" = "0.0₃0868"
    "E20" = "  -3.50%  "
    "E21" = "  -6.33%  "
    "D22" = "64.81"
    "E22" = "  -4.44%  "
    "D23" = "10.25"
    "E23" = "  -3.54%  "
    "D24" = "225.10"
    "E24" = "  -3.90%  "
    "E25" = "  -0.12%  "
    "E26" = "  -6.12%  "
    "D27" = "1.80"
    "E27" = "  +0.26%  "
    "D28" = "22.47"
    "E28" = "  -3.95%  "
    "E29" = "  -1.91%  "
    "E30" = "  -1.64%  "
    "D31" = "149.72"
    "E31" = "  -1.25%  "
    "D32" = "31.45"
    "E32" = "  -7.97%  "
    "E33" = "  -0.15%  "
    "D34" = "4.79"
    "E34" = "  -6.31%  "
    "B35" = "WEMIXToken"
    "C35" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D35" = "2.33"
    "E35" = "  -3.75%  "
    "B36" = "Hedera"
    "C36" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D36" = "0.0690"
    "E36" = "  -4.60%  "
    "E37" = "  -3.87%  "
    "B38" = "Celestia"
    "C38" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D38" = "15.27"
    "E38" = "  -2.37%  "
    "B39" = "Kaspa"
    "C39" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D39" = "0.0960"
    "E39" = "  -2.82%  "
    "E40" = "  -5.46%  "
    "D41" = "1.62"
    "E41" = "  -3.92%  "
    "E42" = "  -5.50%  "
    "D43" = "1.887.92"
    "E43" = "  -3.23%  "
    "E44" = "  -11.74%  "
    "E45" = "  -2.39%  "
    "D46" = "16.07"
    "E46" = "  -9.61%  "
    "D47" = "8.92"
    "E47" = "  -4.54%  "
    "E48" = "  -1.71%  "
    "D49" = "2.400.36"
    "E49" = "  -7.19%  "
    "D50" = "70.77"
    "E50" = "  -0.03%  "
    "D51" = "86.58"
    "E51" = "  -6.10%  "
}

# A new cell value that reads as a plain number (e.g. "81.44") must have
# its cell pre-formatted as Text, otherwise Excel auto-converts it to a
# real number on assignment -- losing the original text formatting and,
# for values such as "1.80", even a significant trailing zero.
function Test-LooksNumeric($text) {
    return [bool]($text -match '^\s*[-+]?[0-9]+(\.[0-9]+)?\s*$')
}

foreach ($address in $updates.Keys) {
    $newValue = $updates[$address]
    $target = $ws.Range($address)
    if (Test-LooksNumeric $newValue) {
        $target.NumberFormat = "@"
    }
    $target.Value = $newValue
}
